# Actualización automática del tracker
# Fill in the "resultado" (G) and "profit" (H) columns for rows whose
# bet outcome has now been settled, and normalize A283 to a real number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 232; Resultado = "Acierto"; Profit = 0.91 },
    @{ Row = 243; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 244; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 245; Resultado = "Acierto"; Profit = 1.63 },
    @{ Row = 250; Resultado = "Acierto"; Profit = 0.91 },
    @{ Row = 251; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 252; Resultado = "Acierto"; Profit = 1.63 },
    @{ Row = 255; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 256; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 257; Resultado = "Acierto"; Profit = 1.75 },
    @{ Row = 258; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 259; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 260; Resultado = "Acierto"; Profit = 2.4 },
    @{ Row = 261; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 262; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 266; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 267; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 268; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 269; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 270; Resultado = "Acierto"; Profit = 1.25 },
    @{ Row = 271; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 274; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 279; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 281; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 282; Resultado = "Fallo";   Profit = -1 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 7).Value = $u.Resultado
    $ws.Cells.Item($u.Row, 8).Value = $u.Profit
}

# Row 283's event_id was stored as text; normalize it back to a number.
$ws.Cells.Item(283, 1).Value = 14452703
